$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched cells in columns B-E so Excel
# doesn't auto-coerce numeric-looking strings (e.g. '2.80' -> 2.8).
$cellList = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "E38", "D39", "E39", "E40", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E46", "E47", "B48", "C48", "D48", "E48", "B49", "C49", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $cellList) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '71.884.63'
$ws.Range("E2").Value = '  +3.19%  '
$ws.Range("D3").Value = '4.043.49'
$ws.Range("E3").Value = '  +2.75%  '
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '523.72'
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").Value = '149.21'
$ws.Range("E6").Value = '  +2.67%  '
$ws.Range("E7").Value = '  +1.21%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").Value = '0.742'
$ws.Range("E9").Value = '  +1.93%  '
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").Value = '0.0000341'
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '46.67'
$ws.Range("E12").Value = '  +9.24%  '
$ws.Range("D13").Value = '10.79'
$ws.Range("E13").Value = '  +3.67%  '
$ws.Range("D14").Value = '4.689.34'
$ws.Range("E14").Value = '  +2.85%  '
$ws.Range("D15").Value = '4.044.35'
$ws.Range("E15").Value = '  +2.74%  '
$ws.Range("D16").Value = '21.51'
$ws.Range("E16").Value = '  +8.32%  '
$ws.Range("D17").Value = '14.32'
$ws.Range("E17").Value = '  +2.01%  '
$ws.Range("D18").Value = '1.23'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("E19").Value = '  -1.74%  '
$ws.Range("D20").Value = '72.120.64'
$ws.Range("E20").Value = '  +3.75%  '
$ws.Range("D21").Value = '445.34'
$ws.Range("E21").Value = '  +2.92%  '
$ws.Range("D22").Value = '3.57'
$ws.Range("E22").Value = '  +5.46%  '
$ws.Range("D23").Value = '94.86'
$ws.Range("E23").Value = '  +6.99%  '
$ws.Range("D24").Value = '14.33'
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").Value = '12.36'
$ws.Range("E25").Value = '  +4.82%  '
$ws.Range("E26").Value = '  -4.10%  '
$ws.Range("D27").Value = '11.16'
$ws.Range("E27").Value = '  +3.19%  '
$ws.Range("D28").Value = '37.22'
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("E29").Value = '  +2.24%  '
$ws.Range("D30").Value = '699.18'
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("E31").Value = '  +2.88%  '
$ws.Range("D32").Value = '2.92'
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("D33").Value = '6.97'
$ws.Range("E33").Value = '  +14.20%  '
$ws.Range("D34").Value = '67.81'
$ws.Range("E34").Value = '  -4.85%  '
$ws.Range("D35").Value = '0.0₃0908'
$ws.Range("E35").Value = '  +5.67%  '
$ws.Range("D36").Value = '0.442'
$ws.Range("E36").Value = '  -5.69%  '
$ws.Range("D37").Value = '40.86'
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("E38").Value = '  +5.61%  '
$ws.Range("D39").Value = '3.54'
$ws.Range("E39").Value = '  +18.07%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("D42").Value = '0.0489'
$ws.Range("E42").Value = '  +1.27%  '
$ws.Range("E43").Value = '  +0.63%  '
$ws.Range("D44").Value = '2.80'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '3.53'
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("E47").Value = '  +2.64%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '9.19'
$ws.Range("E48").Value = '  +5.87%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").Value = '0.000278'
$ws.Range("E49").Value = '  +19.00%  '
$ws.Range("D50").Value = '3.37'
$ws.Range("E50").Value = '  +1.30%  '
$ws.Range("D51").Value = '0.0₆0345'
$ws.Range("E51").Value = '  -1.98%  '
